$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: drop F2, keep D2/E2 text (unchanged text, just reshuffled shared-string slots) ---
$ws.Range("D2").Value = "1 to 100"
$ws.Range("E2").Value = "1 to 10"
$ws.Range("F2").Clear()

# --- Row 3 header: F3 text changes ---
$ws.Range("F3").Value = "Risk response strategy"

# --- Row 4: existing risk row gets reworked with new columns ---
$ws.Range("B4").Value = "Not enought time"
$ws.Range("C4").ClearContents()
$ws.Range("C4").WrapText = $true
$ws.Range("D4").NumberFormat = "0%"
$ws.Range("D4").Value = 0.4
$ws.Range("E4").Value = 10
$ws.Range("F4").ClearContents()
$ws.Range("F4").WrapText = $true

# --- Row 5: new risk row ---
$ws.Range("B5").Value = "Underestimated complexity"
$ws.Range("C5").WrapText = $true
$ws.Range("D5").NumberFormat = "0%"
$ws.Range("D5").Value = 0.2
$ws.Range("E5").Value = 5
$ws.Range("F5").WrapText = $true

# --- Row 6: new risk row with descriptive text in C/F (wraps -> taller row) ---
$ws.Range("B6").Value = "Unreliability of a team member"
$ws.Range("C6").Value = "A team member doesnt complete his work."
$ws.Range("C6").WrapText = $true
$ws.Range("D6").NumberFormat = "0%"
$ws.Range("D6").Value = 0.15
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = "The other team member have to plit up the undone work and complete it."
$ws.Range("F6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 43.5

# --- Rows 7-10: trailing styled-but-empty cells ---
$ws.Range("C7").WrapText = $true
$ws.Range("F7").WrapText = $true
$ws.Range("C8").WrapText = $true
$ws.Range("F8").WrapText = $true
$ws.Range("C9").WrapText = $true
$ws.Range("C10").WrapText = $true

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 26.166666666666668
$ws.Columns.Item(3).ColumnWidth = 27.346354166666668
$ws.Columns.Item(6).ColumnWidth = 28.983072916666668

# --- Selection ---
$ws.Range("C6").Select() | Out-Null
